# Add the "NA" values under the duplicate_image_filename column (column E)
# for every practice/trial row (rows 2-21) of the stimuli sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2:E21").Value = "NA"
